# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# Updates the recalculated term-structure values (columns D:H) that resulted from
# excluding a < $5 option price data point from the calibration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D5" = 113573.6936913498;  "E5" = -0.02356514351686368; "F5" = 0.2274410172975442;  "G5" = -1.23096294590669;   "H5" = 9.76512849579626
    "D6" = 114830.5538012689;  "E6" = -0.02077481980262083; "F6" = 0.2239908128572526;  "G6" = -0.8985853054275793; "H6" = 6.501530830932288
    "D7" = 115416.1813781023;  "E7" = -0.02258236139945237; "F7" = 0.2104007138331737;  "G7" = -0.3540548661700166; "H7" = 4.384467486549173
    "D8" = 115287.7545103493;  "E8" = -0.04060636130763168; "F8" = 0.2166517836764271;  "G8" = -1.603430736963122;  "H8" = 11.25603432029452
    "D9" = 116821.5235108157;  "E9" = -0.06731304907521063; "F9" = 0.3167145157044435;  "G9" = -1.664623964732298;  "H9" = 10.88502509680619
    "D10" = 118247.5070792605; "E10" = -0.11316719516027;   "F10" = 0.4421140365973167; "G10" = -1.872113354886799; "H10" = 9.381189237443115
    "D17" = 111398.9421081011; "E17" = -0.001876497208298536; "F17" = 0.159925512416507; "G17" = -0.6340009709874348; "H17" = 4.74796040866665
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
